# Auto-generated edit script
# Applies updated Leve profit/price figures (scheduled market-data refresh)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 689.6667
$ws.Range("I28").Value = 601.1429000000001
$ws.Range("J28").Value = 999.5
$ws.Range("K28").Value = 601.1429000000001
$ws.Range("L28").Value = 999.5
$ws.Range("M28").Value = -116.1429000000001
$ws.Range("N28").Value = -1969.5
$ws.Range("H39").Value = 202.2
$ws.Range("I39").Value = 178.58824
$ws.Range("K39").Value = 535.76472
$ws.Range("M39").Value = -239.76472
$ws.Range("H116").Value = 6008.7856
$ws.Range("I116").Value = 4374.375
$ws.Range("K116").Value = 4374.375
$ws.Range("M116").Value = -932.375
$ws.Range("H132").Value = 28277.732
$ws.Range("I132").Value = 33488.777
$ws.Range("K132").Value = 100466.331
$ws.Range("M132").Value = -97936.33100000001
$ws.Range("H137").Value = 15544.667
$ws.Range("I137").Value = 1518.7778
$ws.Range("J137").Value = 36583.5
$ws.Range("K137").Value = 4556.3334
$ws.Range("L137").Value = 109750.5
$ws.Range("M137").Value = -2006.3334
$ws.Range("N137").Value = -114850.5

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2341.6667
$ws.Range("I45").Value = 2384.375
$ws.Range("K45").Value = 2384.375
$ws.Range("M45").Value = -2007.375
$ws.Range("H74").Value = 15453.139
$ws.Range("I74").Value = 1694.1923
$ws.Range("J74").Value = 51226.4
$ws.Range("K74").Value = 1694.1923
$ws.Range("L74").Value = 51226.4
$ws.Range("M74").Value = -820.1922999999999
$ws.Range("N74").Value = -52974.4
$ws.Range("H77").Value = 15453.139
$ws.Range("I77").Value = 1694.1923
$ws.Range("J77").Value = 51226.4
$ws.Range("K77").Value = 8470.961499999999
$ws.Range("L77").Value = 256132
$ws.Range("M77").Value = -4102.961499999999
$ws.Range("N77").Value = -264868
$ws.Range("H110").Value = 5685810.5
$ws.Range("I110").Value = 6064820
$ws.Range("K110").Value = 6064820
$ws.Range("M110").Value = -6062775
$ws.Range("H132").Value = 2389644
$ws.Range("I132").Value = 1415.129
$ws.Range("J132").Value = 9120107
$ws.Range("K132").Value = 4245.387
$ws.Range("L132").Value = 27360321
$ws.Range("M132").Value = -1715.387
$ws.Range("N132").Value = -27365381

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 46800
$ws.Range("J130").Value = 46800
$ws.Range("L130").Value = 46800
$ws.Range("N130").Value = -56840

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 13521.833
$ws.Range("I58").Value = 5855.8
$ws.Range("K58").Value = 5855.8
$ws.Range("M58").Value = -5652.8
$ws.Range("H132").Value = 31252664
$ws.Range("I132").Value = 2811.8076
$ws.Range("J132").Value = 166668690
$ws.Range("K132").Value = 8435.4228
$ws.Range("L132").Value = 500006070
$ws.Range("M132").Value = -5905.4228
$ws.Range("N132").Value = -500011130
$ws.Range("H134").Value = 25006554
$ws.Range("I134").Value = 3485
$ws.Range("K134").Value = 10455
$ws.Range("M134").Value = -7920
$ws.Range("H136").Value = 13521.833
$ws.Range("I136").Value = 5855.8
$ws.Range("K136").Value = 17567.4
$ws.Range("M136").Value = -15017.4

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2903.2
$ws.Range("J34").Value = 4748
$ws.Range("L34").Value = 14244
$ws.Range("N34").Value = -14412
$ws.Range("H39").Value = 5869.6
$ws.Range("I39").Value = 3383
$ws.Range("J39").Value = 9599.5
$ws.Range("K39").Value = 10149
$ws.Range("L39").Value = 28798.5
$ws.Range("M39").Value = -9855
$ws.Range("N39").Value = -29386.5
$ws.Range("H55").Value = 3477.7144
$ws.Range("J55").Value = 4481.3335
$ws.Range("L55").Value = 13444.0005
$ws.Range("N55").Value = -13798.0005

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18754.5
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 18754.5
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 18754.5
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -19294.5
$ws.Range("H73").Value = 18754.5
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 18754.5
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 18754.5
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -20626.5
$ws.Range("H113").Value = 2591.2307
$ws.Range("I113").Value = 2048.7778
$ws.Range("K113").Value = 2048.7778
$ws.Range("M113").Value = 121.2222000000002
$ws.Range("H122").Value = 944728.9399999999
$ws.Range("I122").Value = 1133335
$ws.Range("K122").Value = 3400005
$ws.Range("M122").Value = -3397555

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 9061.799999999999
$ws.Range("I32").Value = 1764.6666
$ws.Range("K32").Value = 1764.6666
$ws.Range("M32").Value = -1447.6666
$ws.Range("H43").Value = 502500
$ws.Range("I43").Value = 5000
$ws.Range("K43").Value = 5000
$ws.Range("M43").Value = -4807
$ws.Range("H46").Value = 1754.381
$ws.Range("I46").Value = 730
$ws.Range("J46").Value = 2074.5
$ws.Range("K46").Value = 730
$ws.Range("L46").Value = 2074.5
$ws.Range("M46").Value = -542
$ws.Range("N46").Value = -2450.5
$ws.Range("H55").Value = 1676.56
$ws.Range("I55").Value = 1503.6666
$ws.Range("J55").Value = 1836.1538
$ws.Range("K55").Value = 1503.6666
$ws.Range("L55").Value = 1836.1538
$ws.Range("M55").Value = -1330.6666
$ws.Range("N55").Value = -2182.1538
$ws.Range("H61").Value = 3461.4546
$ws.Range("I61").Value = 3307.6
$ws.Range("K61").Value = 3307.6
$ws.Range("M61").Value = -3105.6
$ws.Range("H113").Value = 3461.4546
$ws.Range("I113").Value = 3307.6
$ws.Range("K113").Value = 3307.6
$ws.Range("M113").Value = -1137.6

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 12373.333
$ws.Range("I62").Value = 9444.444
$ws.Range("J62").Value = 16766.666
$ws.Range("K62").Value = 9444.444
$ws.Range("L62").Value = 16766.666
$ws.Range("M62").Value = -8820.444
$ws.Range("N62").Value = -18014.666
$ws.Range("H65").Value = 12373.333
$ws.Range("I65").Value = 9444.444
$ws.Range("J65").Value = 16766.666
$ws.Range("K65").Value = 47222.22
$ws.Range("L65").Value = 83833.33
$ws.Range("M65").Value = -44102.22
$ws.Range("N65").Value = -90073.33
$ws.Range("H95").Value = 36641.5
$ws.Range("J95").Value = 36641.5
$ws.Range("L95").Value = 36641.5
$ws.Range("N95").Value = -42133.5
$ws.Range("H119").Value = 235000
$ws.Range("J119").Value = 235000
$ws.Range("L119").Value = 235000
$ws.Range("N119").Value = -244676
$ws.Range("H122").Value = 638301.2
$ws.Range("I122").Value = 928062.4399999999
$ws.Range("K122").Value = 2784187.32
$ws.Range("M122").Value = -2781737.32
$ws.Range("H132").Value = 7251
$ws.Range("I132").Value = 2545
$ws.Range("J132").Value = 26075
$ws.Range("K132").Value = 7635
$ws.Range("L132").Value = 78225
$ws.Range("M132").Value = -5105
$ws.Range("N132").Value = -83285
$ws.Range("H136").Value = 14109.223
$ws.Range("I136").Value = 2330.75
$ws.Range("J136").Value = 28832.312
$ws.Range("K136").Value = 6992.25
$ws.Range("L136").Value = 86496.936
$ws.Range("M136").Value = -4442.25
$ws.Range("N136").Value = -91596.936
